$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the formatting of the last existing data row (row 9) down through
# the two new rows (10:11) so the new cells pick up the same styles
# (date format, time format, etc.) already used by the table.
$ws.Range("A9:F9").Copy()
$ws.Range("A10:F11").PasteSpecial(-4122) # xlPasteFormats

# Row 10
$ws.Range("A10").Value = 44149
$ws.Range("B10").Value = 0.64861111111111114
$ws.Range("C10").Value = 0.72013888888888899
$ws.Range("D10").Value = 0
$ws.Range("F10").Value = "coding"

# Row 11
$ws.Range("A11").Value = 44150
$ws.Range("B11").Value = 0.5805555555555556
$ws.Range("C11").Value = 0.60902777777777783
$ws.Range("D11").Value = 0
$ws.Range("F11").Value = "coding"

# Extend the shared formula in column E down through the new rows, same as
# dragging the fill handle on E9 down to E11 in the UI.
$ws.Range("E10").Formula = "=C10-B10-TIME(0,D10,0)"
$ws.Range("E11").Formula = "=C11-B11-TIME(0,D11,0)"

$ws.Range("F12").Select()
